# Removed shlogin test cases from Iam API
#
# The sheet has two obsolete Account-API test rows (OPQA-3561 "redirect to
# password reset page" and OPQA-3562 "redirect to Shibboleth login page" /
# "/account/shlogin") that need to disappear. Deleting the rows shifts every
# row below them up by two and lets Excel renumber the shared-string table
# and dimension/selection bookkeeping on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 53 = OPQA-3561 (account/resetpassword), Row 54 = OPQA-3562 (account/shlogin)
$ws.Rows("53:54").Delete()

# Restore the view state saved with the workbook (scrolled near the bottom,
# cursor sitting on the last data row's B column).
$aw = $excel.ActiveWindow
$aw.ScrollRow = 46
$aw.ScrollColumn = 1
$ws.Range("B56").Select()
